$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Runmode values in column C between rows 2 and 3
$ws.Range("C2").Value = "N"
$ws.Range("C3").Value = "Y"

# Widen column A to match the target width of 27.140625 characters.
# (The host's ColumnWidth setter quantizes to whole pixels, so 26.33 is the
# input that lands closest on the intended stored width.)
$ws.Columns.Item(1).ColumnWidth = 26.333333333333336

# Update the sheet view's active/selected cell to C2
$ws.Range("C2").Select()
